$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.868.56"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.28"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.62"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5078"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2573"
$ws.Range("E8").Value = "  +0.88%  "

$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07759"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.253"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.81"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.855.10"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.75"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7656"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.881.98"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.405"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.96"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.875"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.016"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  +3.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.38"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1243"
$ws.Range("E27").Value = "  +4.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.756"
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04880"
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.238"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.185"
$ws.Range("E33").Value = "  +1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.542"
$ws.Range("E34").Value = "  +1.32%  "

$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8942"
$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5511"
$ws.Range("E37").Value = "  +2.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.540"
$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.117.65"
$ws.Range("E39").Value = "  -1.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01552"
$ws.Range("E40").Value = "  +0.62%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.605"
$ws.Range("E42").Value = "  +2.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7972"
$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.31"
$ws.Range("E44").Value = "  -1.31%  "

$ws.Range("E45").Value = "  -4.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.765.28"
$ws.Range("E46").Value = "  -0.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4445"
$ws.Range("E47").Value = "  -1.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.69"
$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05128"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.553"
$ws.Range("E51").Value = "  +3.13%  "
